$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (46082 -> 46083, i.e. 2026-03-01 -> 2026-03-02) for every data row (2..121).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 1
    }
}
